$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reverser from Data")
$ws.Activate()

# Row 5: change item lookup to "Iron Crown" and overwrite the SUMIFS results
# with the literal values the author typed over them.
$ws.Range("L5").Value = "Iron Crown"

$ws.Range("P5").Value = 0.5
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 40

# Row 6: change item lookup to "Iron Crown" too; its SUMIFS formulas stay intact
# and will recompute automatically against the Data sheet.
$ws.Range("L6").Value = "Iron Crown"

# Recalculate the workbook so all dependent formulas/cached chart values update.
$excel.CalculateFullRebuild()

# Reflect the selection the author left behind after editing L5:U5.
$ws.Range("L5:U5").Select()
$wb.Save()
